# Apply commit: "added mouse feedback and absent trials"
# Adds 100 new trial rows (rows 102-201) to Sheet1 for four new conditions:
# 0_targ_1_dist_red, 0_targ_4_dist_red, 0_targ_10_dist_red, 0_targ_16_dist_red
# Each new row carries tarPos = [0.8,-0.8] and tarSize = [0.2,0.2]

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(102, "0_targ_1_dist_red/trial1.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(103, "0_targ_1_dist_red/trial2.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(104, "0_targ_1_dist_red/trial3.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(105, "0_targ_1_dist_red/trial4.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(106, "0_targ_1_dist_red/trial5.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(107, "0_targ_1_dist_red/trial6.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(108, "0_targ_1_dist_red/trial7.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(109, "0_targ_1_dist_red/trial8.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(110, "0_targ_1_dist_red/trial9.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(111, "0_targ_1_dist_red/trial10.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(112, "0_targ_1_dist_red/trial11.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(113, "0_targ_1_dist_red/trial12.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(114, "0_targ_1_dist_red/trial13.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(115, "0_targ_1_dist_red/trial14.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(116, "0_targ_1_dist_red/trial15.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(117, "0_targ_1_dist_red/trial6.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(118, "0_targ_1_dist_red/trial17.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(119, "0_targ_1_dist_red/trial18.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(120, "0_targ_1_dist_red/trial19.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(121, "0_targ_1_dist_red/trial20.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(122, "0_targ_1_dist_red/trial21.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(123, "0_targ_1_dist_red/trial22.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(124, "0_targ_1_dist_red/trial23.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(125, "0_targ_1_dist_red/trial24.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(126, "0_targ_1_dist_red/trial25.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(127, "0_targ_4_dist_red/trial1.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(128, "0_targ_4_dist_red/trial2.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(129, "0_targ_4_dist_red/trial3.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(130, "0_targ_4_dist_red/trial4.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(131, "0_targ_4_dist_red/trial5.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(132, "0_targ_4_dist_red/trial6.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(133, "0_targ_4_dist_red/trial7.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(134, "0_targ_4_dist_red/trial8.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(135, "0_targ_4_dist_red/trial9.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(136, "0_targ_4_dist_red/trial10.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(137, "0_targ_4_dist_red/trial11.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(138, "0_targ_4_dist_red/trial12.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(139, "0_targ_4_dist_red/trial13.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(140, "0_targ_4_dist_red/trial14.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(141, "0_targ_4_dist_red/trial15.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(142, "0_targ_4_dist_red/trial6.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(143, "0_targ_4_dist_red/trial17.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(144, "0_targ_4_dist_red/trial18.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(145, "0_targ_4_dist_red/trial19.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(146, "0_targ_4_dist_red/trial20.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(147, "0_targ_4_dist_red/trial21.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(148, "0_targ_4_dist_red/trial22.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(149, "0_targ_4_dist_red/trial23.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(150, "0_targ_4_dist_red/trial24.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(151, "0_targ_4_dist_red/trial25.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(152, "0_targ_10_dist_red/trial1.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(153, "0_targ_10_dist_red/trial2.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(154, "0_targ_10_dist_red/trial3.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(155, "0_targ_10_dist_red/trial4.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(156, "0_targ_10_dist_red/trial5.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(157, "0_targ_10_dist_red/trial6.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(158, "0_targ_10_dist_red/trial7.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(159, "0_targ_10_dist_red/trial8.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(160, "0_targ_10_dist_red/trial9.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(161, "0_targ_10_dist_red/trial10.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(162, "0_targ_10_dist_red/trial11.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(163, "0_targ_10_dist_red/trial12.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(164, "0_targ_10_dist_red/trial13.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(165, "0_targ_10_dist_red/trial14.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(166, "0_targ_10_dist_red/trial15.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(167, "0_targ_10_dist_red/trial16.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(168, "0_targ_10_dist_red/trial17.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(169, "0_targ_10_dist_red/trial18.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(170, "0_targ_10_dist_red/trial19.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(171, "0_targ_10_dist_red/trial20.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(172, "0_targ_10_dist_red/trial21.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(173, "0_targ_10_dist_red/trial22.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(174, "0_targ_10_dist_red/trial23.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(175, "0_targ_10_dist_red/trial24.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(176, "0_targ_10_dist_red/trial25.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(177, "0_targ_16_dist_red/trial1.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(178, "0_targ_16_dist_red/trial2.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(179, "0_targ_16_dist_red/trial3.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(180, "0_targ_16_dist_red/trial4.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(181, "0_targ_16_dist_red/trial5.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(182, "0_targ_16_dist_red/trial6.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(183, "0_targ_16_dist_red/trial7.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(184, "0_targ_16_dist_red/trial8.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(185, "0_targ_16_dist_red/trial9.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(186, "0_targ_16_dist_red/trial10.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(187, "0_targ_16_dist_red/trial11.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(188, "0_targ_16_dist_red/trial12.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(189, "0_targ_16_dist_red/trial13.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(190, "0_targ_16_dist_red/trial14.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(191, "0_targ_16_dist_red/trial15.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(192, "0_targ_16_dist_red/trial16.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(193, "0_targ_16_dist_red/trial17.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(194, "0_targ_16_dist_red/trial18.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(195, "0_targ_16_dist_red/trial19.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(196, "0_targ_16_dist_red/trial20.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(197, "0_targ_16_dist_red/trial21.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(198, "0_targ_16_dist_red/trial22.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(199, "0_targ_16_dist_red/trial23.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(200, "0_targ_16_dist_red/trial24.png", "[0.8,-0.8]", "[0.2,0.2]"),
    @(201, "0_targ_16_dist_red/trial25.png", "[0.8,-0.8]", "[0.2,0.2]")
)

# The original author entered the tarPos/tarSize values for the new block
# before the image filenames (e.g. pasted/filled column B & C first), so the
# shared-string table registers "[0.8,-0.8]" and "[0.2,0.2]" right after the
# pre-existing "[9,3]" entry, ahead of the new filenames. Reproduce that by
# writing the first new row's B/C cells before any of the A (filename) cells.
$ws.Cells.Item(102, 2).Value = "[0.8,-0.8]"
$ws.Cells.Item(102, 3).Value = "[0.2,0.2]"

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}

# Update the view so the newly-added rows are visible/selected, mirroring
# the author's on-screen state when the workbook was saved.
$ws.Range("C151:C201").Select()
$excel.ActiveWindow.ScrollRow = 178
$excel.ActiveWindow.ScrollColumn = 1
